$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.002.08"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "2.634.82"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.00%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "597.80"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "153.88"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("D9").Value = "2.633.06"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("E10").Value = "  +10.38%  "
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("E13").Value = "  -0.24%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "27.69"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("E15").Value = "  +3.97%  "
$ws.Range("D16").Value = "3.114.78"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "67.886.17"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "2.632.74"
$ws.Range("E18").Value = "  +0.08%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "375.17"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.39%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "11.38"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.20%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "7.48"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("E23").Value = "  -1.26%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.04"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.44%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "72.37"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.21%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -0.15%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.0000104"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.72%  "
$ws.Range("D29").Value = "2.757.62"
$ws.Range("E30").Value = "  -0.13%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "576.54"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("E32").Value = "  +0.90%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "7.86"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("E34").Value = "  -0.27%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("E37").Value = "  +0.17%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "158.38"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.84%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "19.16"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +5.79%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.369"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("E43").Value = "  +3.73%  "
$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").Value = "0.0₆0321"
$ws.Range("E44").Value = "  +13.34%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "17.13"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +4.83%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "40.46"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.82%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "155.71"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("E50").Value = "  +7.56%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.71"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.45%  "
